# Applies crypto price/volume updates to Sheet1 (rows 2-51),
# matching the daily data refresh performed by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to be stored as plain
# text (mirrors the workbook's existing inline-string cells). Without this,
# Excel auto-detects numeric-looking strings (e.g. "212.90") and silently
# converts them to numbers, which would drop significant trailing zeros.
function Set-TextValue($cellAddress, $text) {
    $range = $ws.Range($cellAddress)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Row 2
Set-TextValue "D2" "26.725.85"
$ws.Range("E2").Value = "  +1.94%  "

# Row 3
Set-TextValue "D3" "1.638.10"
$ws.Range("E3").Value = "  +2.03%  "

# Row 4
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
Set-TextValue "D5" "212.90"
$ws.Range("E5").Value = "  +0.16%  "

# Row 6
Set-TextValue "D6" "0.494"
$ws.Range("E6").Value = "  +1.91%  "

# Row 7
$ws.Range("E7").Value = "  -0.17%  "

# Row 8
Set-TextValue "D8" "0.252"
$ws.Range("E8").Value = "  +1.44%  "

# Row 9
Set-TextValue "D9" "0.0625"
$ws.Range("E9").Value = "  +1.99%  "

# Row 10
Set-TextValue "D10" "19.04"
$ws.Range("E10").Value = "  +4.89%  "

# Row 11
$ws.Range("E11").Value = "  +2.83%  "

# Row 12
Set-TextValue "D12" "1.866.37"
$ws.Range("E12").Value = "  +1.99%  "

# Row 13
Set-TextValue "D13" "1.658.37"
$ws.Range("E13").Value = "  +3.58%  "

# Row 14
Set-TextValue "D14" "4.07"
$ws.Range("E14").Value = "  +1.37%  "

# Row 15
Set-TextValue "D15" "0.525"
$ws.Range("E15").Value = "  +2.49%  "

# Row 16
Set-TextValue "D16" "26.729.42"
$ws.Range("E16").Value = "  +1.92%  "

# Row 17
Set-TextValue "D17" "63.10"
$ws.Range("E17").Value = "  +2.05%  "

# Row 18
Set-TextValue "D18" "0.0₃0742"
$ws.Range("E18").Value = "  +1.93%  "

# Row 19
Set-TextValue "D19" "209.39"
$ws.Range("E19").Value = "  +4.37%  "

# Row 20
$ws.Range("E20").Value = "  -0.22%  "

# Row 21
$ws.Range("E21").Value = "  +1.08%  "

# Row 22
$ws.Range("E22").Value = "  +1.48%  "

# Row 23
Set-TextValue "D23" "6.15"
$ws.Range("E23").Value = "  +2.34%  "

# Row 24
Set-TextValue "D24" "1.95"
$ws.Range("E24").Value = "  +3.97%  "

# Row 25
Set-TextValue "D25" "146.51"
$ws.Range("E25").Value = "  +1.44%  "

# Row 26
$ws.Range("E26").Value = "  -0.18%  "

# Row 27
Set-TextValue "D27" "0.121"
$ws.Range("E27").Value = "  -0.51%  "

# Row 28
Set-TextValue "D28" "6.76"
$ws.Range("E28").Value = "  +3.28%  "

# Row 29
Set-TextValue "D29" "15.41"
$ws.Range("E29").Value = "  +1.52%  "

# Row 30
$ws.Range("E30").Value = "  +5.49%  "

# Row 31
$ws.Range("E31").Value = "  -0.28%  "

# Row 32
$ws.Range("E32").Value = "  +1.03%  "

# Row 33
$ws.Range("E33").Value = "  +1.30%  "

# Row 34
$ws.Range("E34").Value = "  +0.93%  "

# Row 35
$ws.Range("E35").Value = "  +0.21%  "

# Row 36
Set-TextValue "D36" "1.169.34"
$ws.Range("E36").Value = "  +0.62%  "

# Row 37
$ws.Range("E37").Value = "  -0.35%  "

# Row 38
Set-TextValue "D38" "0.810"
$ws.Range("E38").Value = "  +3.03%  "

# Row 39
$ws.Range("E39").Value = "  -0.16%  "

# Row 40
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D40" "0.504"
$ws.Range("E40").Value = "  +1.39%  "

# Row 41
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D41" "2.32"
$ws.Range("E41").Value = "  +0.11%  "

# Row 42
Set-TextValue "D42" "0.796"
$ws.Range("E42").Value = "  +1.29%  "

# Row 43
$ws.Range("E43").Value = "  +1.61%  "

# Row 44
Set-TextValue "D44" "1.775.77"
$ws.Range("E44").Value = "  +1.99%  "

# Row 45
Set-TextValue "D45" "92.45"
$ws.Range("E45").Value = "  +0.79%  "

# Row 46
$ws.Range("E46").Value = "  +2.43%  "

# Row 47
Set-TextValue "D47" "0.0₆0104"
$ws.Range("E47").Value = "  +9.21%  "

# Row 48
Set-TextValue "D48" "54.68"
$ws.Range("E48").Value = "  +1.02%  "

# Row 49
$ws.Range("E49").Value = "  +1.43%  "

# Row 50
$ws.Range("E50").Value = "  +0.38%  "

# Row 51
$ws.Range("E51").Value = "  +4.25%  "
